$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 111-112; existing rows 111+ shift down to 113+
$ws.Rows("111:112").Insert()

# Row 111 - new record
$ws.Cells.Item(111, 1).Value = 8
$ws.Cells.Item(111, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(111, 3).Value = "Coquimbo"
$ws.Cells.Item(111, 4).Value = 44529
$ws.Cells.Item(111, 5).Value = 4
$ws.Cells.Item(111, 6).Value = 100112003
$ws.Cells.Item(111, 7).Value = "Ajo"
$ws.Cells.Item(111, 8).Value = "Chino"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 500
$ws.Cells.Item(111, 11).Value = 19000
$ws.Cells.Item(111, 12).Value = 20000
$ws.Cells.Item(111, 13).Value = 19500
$ws.Cells.Item(111, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(111, 15).Value = "China"
$ws.Cells.Item(111, 16).Value = 1950
$ws.Cells.Item(111, 17).Value = 10
$ws.Cells.Item(111, 18).Value = "Hortaliza"

# Row 112 - new record
$ws.Cells.Item(112, 1).Value = 8
$ws.Cells.Item(112, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(112, 3).Value = "Coquimbo"
$ws.Cells.Item(112, 4).Value = 44529
$ws.Cells.Item(112, 5).Value = 4
$ws.Cells.Item(112, 6).Value = 100112003
$ws.Cells.Item(112, 7).Value = "Ajo"
$ws.Cells.Item(112, 8).Value = "Chino"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 400
$ws.Cells.Item(112, 11).Value = 20000
$ws.Cells.Item(112, 12).Value = 21000
$ws.Cells.Item(112, 13).Value = 20500
$ws.Cells.Item(112, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(112, 15).Value = "China"
$ws.Cells.Item(112, 16).Value = 2050
$ws.Cells.Item(112, 17).Value = 10
$ws.Cells.Item(112, 18).Value = "Hortaliza"
